$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 'мелочный товар'
    3 = 'шелковый товар'
    4 = 'съестной припасы'
    5 = 'бумажный товар'
    6 = 'деревянный товар'
    7 = 'крестьянский товар'
    8 = 'москательный товар'
    9 = 'мелкий товар'
    10 = 'лавочный товар'
    11 = 'рукоделие'
    12 = 'шерстяной товар'
    13 = 'красный товар'
    14 = 'гарусный товар'
    15 = 'особливый товар'
    16 = 'мелочь'
    17 = 'деревенский товар'
    18 = 'серебреный товар'
    19 = 'крамными товар'
    20 = 'небогатый товар'
    21 = 'мясо'
    22 = 'железный товар'
    23 = 'приуготовлять'
    24 = 'щепетильный товар'
    25 = 'пушной товар'
    26 = 'нужный товар'
    27 = 'набойчатый товар'
    28 = 'медный товар'
    29 = 'недорогой товар'
    30 = 'внутренний товар'
    31 = 'питейный припасы'
    32 = 'суровский товар'
    33 = 'оловянный товар'
    34 = 'привозный товар'
    35 = 'купецкий товар'
    36 = 'заморский товар'
    37 = 'произрастание'
    38 = 'галантерейный товар'
    39 = 'надлежащий товар'
    40 = 'харчевой припасы'
    41 = 'меховой товар'
    42 = 'рукодельный товар'
    43 = 'домовый товар'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

